$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 48; this pushes the previous rows 48-162 down to 49-163
# (matching the dimension change from A1:T162 to A1:T163).
$ws.Rows.Item(48).Insert()

# Columns A,B,C,E,F,G,H,I,J,K,L,T hold values that are identical across every
# data row in this sheet (Mercado/Region/Codreg/Tipo/Producto/.../Calidad/Kg
# per unidad). Copy them down from the row directly below (row 49, which after
# the insert still holds the data that used to be in row 48) instead of
# retyping constants.
$constCols = 1,2,3,5,6,7,8,9,10,11,12,20
foreach ($c in $constCols) {
    $ws.Cells.Item(48, $c).Value2 = $ws.Cells.Item(49, $c).Value2
}

# New record values for the inserted row.
$ws.Cells.Item(48, 4).Value2 = 44498                        # D48 Fecha
$ws.Cells.Item(48, 13).Value2 = 530                         # M48 Volumen
$ws.Cells.Item(48, 14).Value2 = 8000                        # N48 Precio minimo
$ws.Cells.Item(48, 15).Value2 = 9000                        # O48 Precio maximo
$ws.Cells.Item(48, 16).Value2 = 8340                        # P48 Precio promedio ponderado
$ws.Cells.Item(48, 17).Value2 = '$/bandeja 7 kilos'         # Q48 Unidad de comercializacion
$ws.Cells.Item(48, 18).Value2 = 'Provincia de Melipilla'    # R48 Origen
$ws.Cells.Item(48, 19).Value2 = 1191                        # S48 Precio $/Kg

# Match the date-number-format used by the other cells in the Fecha column.
$ws.Cells.Item(48, 4).NumberFormat = $ws.Cells.Item(49, 4).NumberFormat
